$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, centered, bordered) to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row data: row number, I value, J value
$data = @(
    @(2,8,9),
    @(3,7,8),
    @(4,8,8),
    @(5,7,7),
    @(6,9,9),
    @(7,8,8),
    @(8,9,10),
    @(9,8,8),
    @(10,9,9),
    @(11,9,9),
    @(12,8,9),
    @(13,9,9),
    @(14,9,9),
    @(15,9,9),
    @(16,9,9),
    @(17,8,9),
    @(18,9,9),
    @(19,11,12),
    @(20,9,9),
    @(21,9,9),
    @(22,9,9),
    @(23,9,9),
    @(24,10,10),
    @(25,7,7),
    @(26,8,8),
    @(27,9,9),
    @(28,10,10),
    @(29,7,7),
    @(30,7,7),
    @(31,9,9),
    @(32,7,7),
    @(33,8,8),
    @(34,8,8),
    @(35,9,9),
    @(36,7,7),
    @(37,10,11),
    @(38,8,8),
    @(39,5,7),
    @(40,6,7),
    @(41,6,7),
    @(42,8,8),
    @(43,5,6),
    @(44,7,7),
    @(45,7,7),
    @(46,7,7),
    @(47,7,7),
    @(48,6,7),
    @(49,7,7),
    @(50,7,8),
    @(51,3,6),
    @(52,8,8),
    @(53,5,6),
    @(54,6,7),
    @(55,7,7),
    @(56,9,9),
    @(57,9,9),
    @(58,7,8),
    @(59,6,6),
    @(60,7,8),
    @(61,5,5),
    @(62,9,9),
    @(63,7,8),
    @(64,8,9),
    @(65,8,9),
    @(66,8,8),
    @(67,4,4),
    @(68,8,8),
    @(69,8,8),
    @(70,5,5)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Output "Done: I0/IF columns added"
